# Add data organization files for MESS
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("samples_retained")

# Row 20 (LimaCastroScott): add modality, language, and notes info
$ws.Range("C20").Value = "acted"
$ws.Range("K20").Value = "non-speech vocalizations"
$ws.Range("G20").Value = "Portuguese*"

# Row 22 (MAV): add modality, language, and notes info
$ws.Range("C22").Value = "acted"
$ws.Range("K22").Value = "non-speech vocalizations"
$ws.Range("G22").Value = "French*"

# Row 31: add modality, language, and notes info
$ws.Range("C31").Value = "acted"
$ws.Range("K31").Value = "non-speech vocalizations"
$ws.Range("G31").Value = "English*"

# Row 24 (MESS): update sample counts, recoded note
$ws.Range("D24").Value = 809
$ws.Range("F24").Value = 91
$ws.Range("K24").Value = "only 1 Canadian speaker; calm is positive here; recoded 91 calm samples as neutral"

# Update the view state to match final selection
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("D31").Select()
